$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H98").Value = 761.8
$ws.Range("I98").Value = 761.8
$ws.Range("K98").Value = 761.8
$ws.Range("M98").Value = 736.2
$ws.Range("H106").Value = 2066.348
$ws.Range("I106").Value = 1258.6666
$ws.Range("K106").Value = 1258.6666
$ws.Range("M106").Value = -627.6666
$ws.Range("H107").Value = 961.8570999999999
$ws.Range("I107").Value = 971
$ws.Range("K107").Value = 971
$ws.Range("M107").Value = 949
$ws.Range("H111").Value = 2185.5386
$ws.Range("I111").Value = 868.2353000000001
$ws.Range("J111").Value = 4673.778
$ws.Range("K111").Value = 2604.7059
$ws.Range("L111").Value = 14021.334
$ws.Range("M111").Value = 462.2941000000001
$ws.Range("N111").Value = -20155.334
$ws.Range("H122").Value = 761.8
$ws.Range("I122").Value = 761.8
$ws.Range("K122").Value = 2285.4
$ws.Range("M122").Value = 164.6000000000004
$ws.Range("H124").Value = 30780
$ws.Range("J124").Value = 30780
$ws.Range("L124").Value = 30780
$ws.Range("N124").Value = -40600
$ws.Range("H129").Value = 952.678
$ws.Range("J129").Value = 1029.7255
$ws.Range("L129").Value = 3089.1765
$ws.Range("N129").Value = -13089.1765
$ws.Range("H137").Value = 1968.6274
$ws.Range("I137").Value = 2076.0322
$ws.Range("J137").Value = 1802.15
$ws.Range("K137").Value = 6228.096600000001
$ws.Range("L137").Value = 5406.450000000001
$ws.Range("M137").Value = -3678.096600000001
$ws.Range("N137").Value = -10506.45
$ws.Range("H141").Value = 3817.5
$ws.Range("I141").Value = 1908
$ws.Range("K141").Value = 5724
$ws.Range("M141").Value = -544

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 3133.9714
$ws.Range("I61").Value = 2851.1904
$ws.Range("K61").Value = 2851.1904
$ws.Range("M61").Value = -2639.1904
$ws.Range("H97").Value = 1465.3846
$ws.Range("I97").Value = 1513.6364
$ws.Range("J97").Value = 1200
$ws.Range("K97").Value = 1513.6364
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -1017.6364
$ws.Range("N97").Value = -2192
$ws.Range("H110").Value = 1888.6666
$ws.Range("I110").Value = 1785.4286
$ws.Range("J110").Value = 2250
$ws.Range("K110").Value = 1785.4286
$ws.Range("L110").Value = 2250
$ws.Range("M110").Value = 259.5714
$ws.Range("N110").Value = -6340
$ws.Range("H132").Value = 19671.285
$ws.Range("I132").Value = 1359.2
$ws.Range("K132").Value = 4077.6
$ws.Range("M132").Value = -1547.6
$ws.Range("H136").Value = 3133.9714
$ws.Range("I136").Value = 2851.1904
$ws.Range("K136").Value = 8553.5712
$ws.Range("M136").Value = -6003.5712

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H94").Value = 2819.75
$ws.Range("I94").Value = 1255.2222
$ws.Range("J94").Value = 4831.2856
$ws.Range("K94").Value = 1255.2222
$ws.Range("L94").Value = 4831.2856
$ws.Range("M94").Value = -804.2221999999999
$ws.Range("N94").Value = -5733.2856

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H99").Value = 5566.1904
$ws.Range("I99").Value = 4126.364
$ws.Range("J99").Value = 7150
$ws.Range("K99").Value = 4126.364
$ws.Range("L99").Value = 7150
$ws.Range("M99").Value = -2628.364
$ws.Range("N99").Value = -10146
$ws.Range("H126").Value = 5566.1904
$ws.Range("I126").Value = 4126.364
$ws.Range("J126").Value = 7150
$ws.Range("K126").Value = 12379.092
$ws.Range("L126").Value = 21450
$ws.Range("M126").Value = -9909.091999999999
$ws.Range("N126").Value = -26390
$ws.Range("H132").Value = 20651.643
$ws.Range("I132").Value = 42926.082
$ws.Range("J132").Value = 3945.8125
$ws.Range("K132").Value = 128778.246
$ws.Range("L132").Value = 11837.4375
$ws.Range("M132").Value = -126248.246
$ws.Range("N132").Value = -16897.4375

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H5").Value = 849.36
$ws.Range("I5").Value = 993.4286
$ws.Range("J5").Value = 793.3333
$ws.Range("K5").Value = 2980.2858
$ws.Range("L5").Value = 2379.9999
$ws.Range("M5").Value = -2868.2858
$ws.Range("N5").Value = -2603.9999
$ws.Range("H68").Value = 4100.3145
$ws.Range("J68").Value = 5941
$ws.Range("L68").Value = 17823
$ws.Range("N68").Value = -19445
$ws.Range("H71").Value = 4100.3145
$ws.Range("J71").Value = 5941
$ws.Range("L71").Value = 53469
$ws.Range("N71").Value = -61581
$ws.Range("H92").Value = 1057
$ws.Range("I92").Value = 550
$ws.Range("J92").Value = 1259.8
$ws.Range("K92").Value = 1650
$ws.Range("L92").Value = 3779.4
$ws.Range("M92").Value = -402
$ws.Range("N92").Value = -6275.4
$ws.Range("H131").Value = 787.7
$ws.Range("J131").Value = 804.74225
$ws.Range("L131").Value = 2414.22675
$ws.Range("N131").Value = -12494.22675
$ws.Range("H135").Value = 849.36
$ws.Range("I135").Value = 993.4286
$ws.Range("J135").Value = 793.3333
$ws.Range("K135").Value = 8940.857399999999
$ws.Range("L135").Value = 7139.9997
$ws.Range("M135").Value = -6405.857399999999
$ws.Range("N135").Value = -12209.9997

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 2268
$ws.Range("I102").Value = 2216
$ws.Range("J102").Value = 2450
$ws.Range("K102").Value = 2216
$ws.Range("L102").Value = 2450
$ws.Range("M102").Value = -594
$ws.Range("N102").Value = -5694
$ws.Range("H126").Value = 5291.1143
$ws.Range("I126").Value = 5476.316
$ws.Range("J126").Value = 5071.1875
$ws.Range("K126").Value = 16428.948
$ws.Range("L126").Value = 15213.5625
$ws.Range("M126").Value = -13958.948
$ws.Range("N126").Value = -20153.5625

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 4368.077
$ws.Range("I7").Value = 5403.3335
$ws.Range("J7").Value = 2956.3635
$ws.Range("K7").Value = 5403.3335
$ws.Range("L7").Value = 2956.3635
$ws.Range("M7").Value = -5291.3335
$ws.Range("N7").Value = -3180.3635
$ws.Range("H40").Value = 226922.2
$ws.Range("I40").Value = 281402.75
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 281402.75
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -281266.75
$ws.Range("N40").Value = -9272
$ws.Range("H46").Value = 1715.4706
$ws.Range("J46").Value = 1900
$ws.Range("L46").Value = 1900
$ws.Range("N46").Value = -2276
$ws.Range("H126").Value = 4368.077
$ws.Range("I126").Value = 5403.3335
$ws.Range("J126").Value = 2956.3635
$ws.Range("K126").Value = 16210.0005
$ws.Range("L126").Value = 8869.0905
$ws.Range("M126").Value = -13740.0005
$ws.Range("N126").Value = -13809.0905

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H107").Value = 2522.6667
$ws.Range("I107").Value = 1926
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 5778
$ws.Range("L107").Value = 9000
$ws.Range("M107").Value = -3858
$ws.Range("N107").Value = -12840
$ws.Range("H126").Value = 1546.8
$ws.Range("I126").Value = 1139.4
$ws.Range("J126").Value = 1954.2
$ws.Range("K126").Value = 3418.2
$ws.Range("L126").Value = 5862.6
$ws.Range("M126").Value = -948.2000000000003
$ws.Range("N126").Value = -10802.6
$ws.Range("H132").Value = 1949.9667
$ws.Range("I132").Value = 1772.7727
$ws.Range("J132").Value = 2437.25
$ws.Range("K132").Value = 5318.3181
$ws.Range("L132").Value = 7311.75
$ws.Range("M132").Value = -2788.3181
$ws.Range("N132").Value = -12371.75
